# ---------------------------------------------------------------------------
# Edit: Resources/.../Component three/C1--C2-and-C3-PowerPoint.pptx
#
# 1) Slide 16's table (cash-flow "PLENARY" table) switches from the custom
#    "Table_0" style to the built-in "Medium Style 2 - Accent 1" table style.
# 2) The deck's theme colour scheme (the one used by the Slide Master /
#    "Integral" design) is recoloured to the stock "Office Theme" palette.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16, shape 3 (the graphicFrame holding the table)
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{31E4ED0B-4293-4790-A00F-B9105925763E}")

# --- 2) Recolour the theme's colour scheme to the default "Office" palette.
function HexToBgrInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$scheme = $p.SlideMaster.Theme.ThemeColorScheme
$scheme.Item(1).RGB  = HexToBgrInt("000000")  # dk1
$scheme.Item(2).RGB  = HexToBgrInt("FFFFFF")  # lt1
$scheme.Item(3).RGB  = HexToBgrInt("44546A")  # dk2
$scheme.Item(4).RGB  = HexToBgrInt("E7E6E6")  # lt2
$scheme.Item(5).RGB  = HexToBgrInt("5B9BD5")  # accent1
$scheme.Item(6).RGB  = HexToBgrInt("ED7D31")  # accent2
$scheme.Item(7).RGB  = HexToBgrInt("A5A5A5")  # accent3
$scheme.Item(8).RGB  = HexToBgrInt("FFC000")  # accent4
$scheme.Item(9).RGB  = HexToBgrInt("4472C4")  # accent5
$scheme.Item(10).RGB = HexToBgrInt("70AD47")  # accent6
$scheme.Item(11).RGB = HexToBgrInt("0563C1")  # hyperlink
$scheme.Item(12).RGB = HexToBgrInt("954F72")  # followed hyperlink
